# Applies the "output generated at 456a3b4" data refresh to 杭州-漫展信息.xlsx
# Sheets: 1=展览 (exhibitions), 2=演出 (performances), 3=本地生活 (local life, empty), 4=全部类型 (all types, union of the above)

$wb = $excel.ActiveWorkbook

$wsExpo = $wb.Worksheets.Item("展览")
$wsShow = $wb.Worksheets.Item("演出")
$wsAll  = $wb.Worksheets.Item("全部类型")

# ---------------------------------------------------------------------------
# 1) 展览 (sheet 1) — "想去人数" (F column) refreshed counts
# ---------------------------------------------------------------------------
$expoUpdates = @{
    2  = 10593
    4  = 293
    7  = 799
    9  = 1216
    10 = 1138
    11 = 3300
    12 = 2444
    14 = 2229
    15 = 243
    16 = 1941
    18 = 1607
    19 = 604
    20 = 76
    21 = 270
    23 = 33
    24 = 251
    26 = 398
    28 = 82
    29 = 419
    30 = 620
    31 = 51
    32 = 60
    33 = 302
    34 = 21
    35 = 1588
    36 = 595
    37 = 584
    38 = 1783
    39 = 170
    40 = 464
    41 = 66
    42 = 496
    43 = 1099
    45 = 377
}

foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Range("F$row").Value = $expoUpdates[$row]
}

# ---------------------------------------------------------------------------
# 2) 演出 (sheet 2) — row 2 ("天空之城" concert) becomes sold out, row 9 count +1
# ---------------------------------------------------------------------------
$wsShow.Range("G2").Value = "不可售"
$wsShow.Range("F9").Value = 8

# ---------------------------------------------------------------------------
# 3) 全部类型 (sheet 4) — "想去人数" (F column) refreshed counts (shared rows
#    with 展览, offset by the two 演出 rows inserted near the top)
# ---------------------------------------------------------------------------
$allUpdates = @{
    2  = 10593
    5  = 293
    8  = 799
    9  = 1138
    10 = 3300
    11 = 2444
    12 = 2229
    13 = 1941
    15 = 1607
    16 = 604
    17 = 76
    18 = 270
    20 = 33
    21 = 251
    23 = 398
    25 = 82
    26 = 419
    27 = 620
    28 = 51
    32 = 60
    33 = 302
    34 = 21
    35 = 1588
    36 = 595
    38 = 584
    39 = 1783
    40 = 170
    42 = 8
    44 = 464
    45 = 66
    46 = 496
    47 = 1099
    48 = 377
}

foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}

# Row 3 on 全部类型 now holds the data that used to be on row 4 (星玫Rostar idol
# off-event), and row 4 is replaced by a brand-new row pulled in from 演出
# (浙江·燃爆全场 world-film-themed concert), matching 演出!A3 which was already
# present but not yet merged into the "all types" sheet.

$wsAll.Range("C3").Value = "杭州·星玫Rostar偶像团 1st off会 - 莫里生日SP"
$wsAll.Range("D3").Value = "下沙大道30号 杭州璞砚酒店"
$wsAll.Range("E3").Value = "2024.03.23 12:00-03.23 21:00"
$wsAll.Range("F3").Value = 23
$wsAll.Range("G3").Value = 58
$wsAll.Range("H3").Value = "https://show.bilibili.com/platform/detail.html?id=82690"
$wsAll.Range("I3").Value = "//i2.hdslb.com/bfs/openplatform/202403/grcA9fYK1710327154137.jpeg"

$wsAll.Range("C4").Value = "浙江·燃爆全场·世界电影主题音乐会 《复仇者联盟》、《歌剧魅影》、《泰坦尼克号》燃情主题音乐"
$wsAll.Range("D4").Value = "曙光路31号 浙江音乐厅"
$wsAll.Range("E4").Value = "2024.03.23 19:30-03.23 21:00"
$wsAll.Range("F4").Value = 1
$wsAll.Range("G4").Value = "不可售"
$wsAll.Range("H4").Value = "https://show.bilibili.com/platform/detail.html?id=82685"
$wsAll.Range("I4").Value = "//i0.hdslb.com/bfs/openplatform/202403/KFRQDTnB1710210073027.jpeg"
